$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hảo ĐN")

# Row 5 (Dell insprion 3511): record a partial payment of 1,600,000 on 09/09/2023
$ws.Range("E5").Value = 1600000
$ws.Range("F5").Value = 45178

# Row 6 (latitude 5400) is no longer hidden
$ws.Rows.Item(6).Hidden = $false

# Row 7 (Dell latitude 5400) entry is removed entirely
$ws.Range("A7:D7").Clear()
$ws.Range("I7").ClearContents()

# Update the selection to reflect the reviewed range
$ws.Range("A1:J9").Select()
